$wb = $excel.ActiveWorkbook
$raw = $wb.Worksheets.Item("Raw data")

# New sheet goes right after "Raw data", matching the target tab order.
$ws = $wb.Worksheets.Add($null, $raw)
$ws.Name = "Horizontal_Data"

# Match "Raw data"'s outline/page setup so the new sheet looks the same.
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1
$ws.PageSetup.LeftMargin = 54
$ws.PageSetup.RightMargin = 54
$ws.PageSetup.TopMargin = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# Column A is wide (holds the long horizontal rows of data), B:M are narrow.
$ws.Columns.Item(1).ColumnWidth = 80.75
$ws.Range("B1:M1").ColumnWidth = 6.25

# Pull the same cell formatting ("Month"/year label style, highlighted cell)
# that "Raw data" uses, so the new sheet keeps a matching look.
$raw.Range("A1:A5").Copy()
$ws.Range("A1:A5").PasteSpecial(-4122)

$raw.Range("B5").Copy()
$ws.Range("B5").PasteSpecial(-4122)

# Header row: Month, then the 12 month abbreviations.
$ws.Range("A1").Value = "Month"
$ws.Range("B1").Value = " JAN"
$ws.Range("C1").Value = " FEB"
$ws.Range("D1").Value = " MAR"
$ws.Range("E1").Value = " APR"
$ws.Range("F1").Value = " MAY"
$ws.Range("G1").Value = " JUN"
$ws.Range("H1").Value = " JUL"
$ws.Range("I1").Value = " AUG"
$ws.Range("J1").Value = " SEP"
$ws.Range("K1").Value = " OCT"
$ws.Range("L1").Value = " NOV"
$ws.Range("M1").Value = " DEC"

# Each year's twelve monthly values laid out horizontally in a single cell.
$ws.Range("A2").Value = '"1958", 340, 318, 362, 348, 363, 435, 491, 505, 404, 359, 310, 337'
$ws.Range("A3").Value = '"1959", 360, 342, 406, 396, 420, 472, 548, 559, 463, 407, 362, 405'
$ws.Range("A4").Value = '"1960", 417, 391, 419, 461, 472, 535, 622, 606, 508, 461, 390, 432'

$excel.CutCopyMode = $false
$raw.Activate()
